$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price / 1h-volume-change figures (GitHub Actions refresh).
# Each row: row number, new Price (column D, or $null if unchanged), new Volume(1h) text (column E).
$updates = @(
    @{ Row = 2; D = "27.196.33"; E = "  +5.43%  " },
    @{ Row = 3; D = "1.882.01"; E = "  +3.78%  " },
    @{ Row = 4; D = "0.9998"; E = "  -0.01%  " },
    @{ Row = 5; D = "281.41"; E = "  +1.68%  " },
    @{ Row = 6; D = $null; E = "  +0.06%  " },
    @{ Row = 7; D = "0.5298"; E = "  +3.56%  " },
    @{ Row = 8; D = "0.3539"; E = "  +0.50%  " },
    @{ Row = 9; D = "45.53"; E = "  +1.82%  " },
    @{ Row = 10; D = "0.07039"; E = "  +5.48%  " },
    @{ Row = 11; D = "20.40"; E = "  +1.44%  " },
    @{ Row = 12; D = "0.8230"; E = "  -1.42%  " },
    @{ Row = 13; D = "0.07820"; E = "  -0.09%  " },
    @{ Row = 14; D = "1.886.67"; E = "  +4.13%  " },
    @{ Row = 15; D = "90.82"; E = "  +3.25%  " },
    @{ Row = 16; D = "5.204"; E = "  +2.49%  " },
    @{ Row = 17; D = "0.9997"; E = "  +0.07%  " },
    @{ Row = 18; D = "14.63"; E = "  +5.26%  " },
    @{ Row = 19; D = "0.000008177"; E = "  +1.91%  " },
    @{ Row = 20; D = "0.9992"; E = "  +0.00%  " },
    @{ Row = 21; D = "27.223.27"; E = "  +5.23%  " },
    @{ Row = 22; D = "2.133.54"; E = "  +5.13%  " },
    @{ Row = 23; D = "4.774"; E = "  +0.88%  " },
    @{ Row = 24; D = "10.17"; E = "  +1.45%  " },
    @{ Row = 25; D = "6.255"; E = "  +3.15%  " },
    @{ Row = 26; D = "2.407"; E = "  +9.13%  " },
    @{ Row = 27; D = "147.31"; E = "  +4.04%  " },
    @{ Row = 28; D = "17.61"; E = "  +3.34%  " },
    @{ Row = 29; D = "1.676"; E = "  +1.29%  " },
    @{ Row = 30; D = "114.68"; E = "  +5.35%  " },
    @{ Row = 31; D = "4.422"; E = "  +1.77%  " },
    @{ Row = 32; D = "4.389"; E = "  +3.78%  " },
    @{ Row = 33; D = "0.08959"; E = "  +1.96%  " },
    @{ Row = 34; D = "0.04947"; E = "  +1.32%  " },
    @{ Row = 35; D = "1.184"; E = "  +3.75%  " },
    @{ Row = 36; D = "0.7496"; E = "  +2.14%  " },
    @{ Row = 37; D = "2.907"; E = "  +0.50%  " },
    @{ Row = 38; D = "3.315"; E = "  +8.36%  " },
    @{ Row = 39; D = "2.416"; E = "  +5.28%  " },
    @{ Row = 40; D = "0.5318"; E = "  +1.30%  " },
    @{ Row = 41; D = "0.01888"; E = "  +1.61%  " },
    @{ Row = 42; D = "0.9745"; E = "  +2.36%  " },
    @{ Row = 43; D = "117.02"; E = "  +4.47%  " },
    @{ Row = 44; D = "6.322"; E = "  +2.23%  " },
    @{ Row = 45; D = "8.234"; E = "  +1.08%  " },
    @{ Row = 46; D = "0.9997"; E = "  +0.10%  " },
    @{ Row = 47; D = "0.4615"; E = "  +0.76%  " },
    @{ Row = 48; D = $null; E = "  -0.72%  " },
    @{ Row = 49; D = "9.496"; E = "  +2.18%  " },
    @{ Row = 50; D = "36.73"; E = "  +1.43%  " },
    @{ Row = 51; D = "1.532"; E = "  +1.89%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # Force plain text so values like "1.886.67" or "0.9998" are not
        # reinterpreted as numbers/dates by Excel's smart entry.
        $ws.Range("D" + $u.Row).NumberFormat = "@"
        $ws.Range("D" + $u.Row).Value = $u.D
    }
    $ws.Range("E" + $u.Row).Value = $u.E
}
